$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade record as row 8, mirroring the existing rows' layout/styles
# (date column formatted as a date, IsShortSell formatted/typed like the others).
# Copying row 7 onto row 8 first carries over the per-column styles (date
# number format on A, boolean style on G, etc.) exactly like the other data
# rows, then we overwrite the copied values with the new trade's data.
$ws.Range("A7:I7").Copy($ws.Range("A8:I8"))

$ws.Range("A8").Value = 42650.366944444446
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = 10242.129999999999
$ws.Range("D8").Value = 10321.61
$ws.Range("E8").Value = 104.839996
$ws.Range("F8").Value = 104.029999
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = -0.77
$ws.Range("I8").Value = $true

# The extra row's longer BuyPrice/SellPrice numbers widen columns E:F.
$ws.Columns("E:F").AutoFit()
